$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04936366666666667
$ws.Range("H2").Value = 0.148091
$ws.Range("I2").Value = 0.04616170608573571
$ws.Range("J2").Value = 0.0461617060857357
$ws.Range("M2").Value = 8.906151666666668
$ws.Range("N2").Value = 26.718455
$ws.Range("O2").Value = 0.1245005002255258
$ws.Range("P2").Value = 0.1245005002255258
$ws.Range("Q2").Value = 0.4396403021561112
$ws.Range("R2").Value = 3.956762719405
$ws.Range("S2").Value = 0.005747155498937793
$ws.Range("T2").Value = 0.005747155498937792
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04936366666666667
$ws.Range("H3").Value = 0.148091
$ws.Range("I3").Value = 0.04616170608573571
$ws.Range("J3").Value = 0.0461617060857357
$ws.Range("M3").Value = 0.4515893333333333
$ws.Range("O3").Value = 0.006312838586270617
$ws.Range("P3").Value = 0.006312838586270617
$ws.Range("Q3").Value = 0.02229210532088889
$ws.Range("R3").Value = 0.200628947888
$ws.Range("S3").Value = 0.0002914113993861155
$ws.Range("T3").Value = 0.0002914113993861155
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04936366666666667
$ws.Range("H4").Value = 0.148091
$ws.Range("I4").Value = 0.04616170608573571
$ws.Range("J4").Value = 0.0461617060857357
$ws.Range("M4").Value = 7.781650666666667
$ws.Range("N4").Value = 23.344952
$ws.Range("O4").Value = 0.1087809232135948
$ws.Range("P4").Value = 0.1087809232135948
$ws.Range("Q4").Value = 0.3841308096257778
$ws.Range("R4").Value = 3.457177286632
$ws.Range("S4").Value = 0.005021513005120948
$ws.Range("T4").Value = 0.005021513005120947
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04936366666666667
$ws.Range("H5").Value = 0.148091
$ws.Range("I5").Value = 0.04616170608573571
$ws.Range("J5").Value = 0.0461617060857357
$ws.Range("M5").Value = 3.892567333333333
$ws.Range("N5").Value = 11.677702
$ws.Range("O5").Value = 0.05441481329981927
$ws.Range("P5").Value = 0.05441481329981927
$ws.Range("Q5").Value = 0.1921513963202222
$ws.Range("R5").Value = 1.729362566882
$ws.Range("S5").Value = 0.00251188061825644
$ws.Range("T5").Value = 0.002511880618256439
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04936366666666667
$ws.Range("H6").Value = 0.148091
$ws.Range("I6").Value = 0.04616170608573571
$ws.Range("J6").Value = 0.0461617060857357
$ws.Range("M6").Value = 43.49559133333333
$ws.Range("N6").Value = 130.486774
$ws.Range("O6").Value = 0.6080317382054886
$ws.Range("P6").Value = 0.6080317382054886
$ws.Range("Q6").Value = 2.147101872048222
$ws.Range("R6").Value = 19.323916848434
$ws.Range("S6").Value = 0.02806778238984076
$ws.Range("T6").Value = 0.02806778238984076
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04936366666666667
$ws.Range("H7").Value = 0.148091
$ws.Range("I7").Value = 0.04616170608573571
$ws.Range("J7").Value = 0.0461617060857357
$ws.Range("M7").Value = 7.007517
$ws.Range("N7").Value = 21.022551
$ws.Range("O7").Value = 0.09795918646930096
$ws.Range("P7").Value = 0.09795918646930096
$ws.Range("Q7").Value = 0.345916733349
$ws.Range("R7").Value = 3.113250600141
$ws.Range("S7").Value = 0.004521963174193649
$ws.Range("T7").Value = 0.004521963174193649
$ws.Range("I8").Value = 0.008057748967298944
$ws.Range("J8").Value = 0.008057748967298944
$ws.Range("M8").Value = 8.906151666666668
$ws.Range("N8").Value = 26.718455
$ws.Range("O8").Value = 0.1245005002255258
$ws.Range("P8").Value = 0.1245005002255258
$ws.Range("Q8").Value = 0.07674134019444445
$ws.Range("R8").Value = 0.6906720617500001
$ws.Range("S8").Value = 0.001003193777120432
$ws.Range("T8").Value = 0.001003193777120432
$ws.Range("I9").Value = 0.008057748967298944
$ws.Range("J9").Value = 0.008057748967298944
$ws.Range("M9").Value = 0.4515893333333333
$ws.Range("O9").Value = 0.006312838586270617
$ws.Range("P9").Value = 0.006312838586270617
$ws.Range("R9").Value = 0.0350207528
$ws.Range("S9").Value = 0.00005086726859924699
$ws.Range("T9").Value = 0.00005086726859924699
$ws.Range("I10").Value = 0.008057748967298944
$ws.Range("J10").Value = 0.008057748967298944
$ws.Range("M10").Value = 7.781650666666667
$ws.Range("N10").Value = 23.344952
$ws.Range("O10").Value = 0.1087809232135948
$ws.Range("P10").Value = 0.1087809232135948
$ws.Range("Q10").Value = 0.06705188991111111
$ws.Range("R10").Value = 0.6034670092
$ws.Range("S10").Value = 0.0008765293716861692
$ws.Range("T10").Value = 0.0008765293716861692
$ws.Range("I11").Value = 0.008057748967298944
$ws.Range("J11").Value = 0.008057748967298944
$ws.Range("M11").Value = 3.892567333333333
$ws.Range("N11").Value = 11.677702
$ws.Range("O11").Value = 0.05441481329981927
$ws.Range("P11").Value = 0.05441481329981927
$ws.Range("Q11").Value = 0.03354095518888889
$ws.Range("R11").Value = 0.3018685967
$ws.Range("S11").Value = 0.0004384609056723836
$ws.Range("T11").Value = 0.0004384609056723836
$ws.Range("I12").Value = 0.008057748967298944
$ws.Range("J12").Value = 0.008057748967298944
$ws.Range("M12").Value = 43.49559133333333
$ws.Range("N12").Value = 130.486774
$ws.Range("O12").Value = 0.6080317382054886
$ws.Range("P12").Value = 0.6080317382054886
$ws.Range("Q12").Value = 0.3747870119888889
$ws.Range("R12").Value = 3.3730831079
$ws.Range("S12").Value = 0.004899367110610258
$ws.Range("T12").Value = 0.004899367110610258
$ws.Range("I13").Value = 0.008057748967298944
$ws.Range("J13").Value = 0.008057748967298944
$ws.Range("M13").Value = 7.007517
$ws.Range("N13").Value = 21.022551
$ws.Range("O13").Value = 0.09795918646930096
$ws.Range("P13").Value = 0.09795918646930096
$ws.Range("Q13").Value = 0.06038143815
$ws.Range("R13").Value = 0.54343294335
$ws.Range("S13").Value = 0.0007893305336104545
$ws.Range("T13").Value = 0.0007893305336104545
$ws.Range("G14").Value = 1.011383666666666
$ws.Range("H14").Value = 3.034151
$ws.Range("I14").Value = 0.9457805449469654
$ws.Range("J14").Value = 0.9457805449469653
$ws.Range("M14").Value = 8.906151666666668
$ws.Range("N14").Value = 26.718455
$ws.Range("O14").Value = 0.1245005002255258
$ws.Range("P14").Value = 0.1245005002255258
$ws.Range("Q14").Value = 9.007536328522777
$ws.Range("R14").Value = 81.06782695670499
$ws.Range("S14").Value = 0.1177501509494676
$ws.Range("T14").Value = 0.1177501509494675
$ws.Range("G15").Value = 1.011383666666666
$ws.Range("H15").Value = 3.034151
$ws.Range("I15").Value = 0.9457805449469654
$ws.Range("J15").Value = 0.9457805449469653
$ws.Range("M15").Value = 0.4515893333333333
$ws.Range("O15").Value = 0.006312838586270617
$ws.Range("P15").Value = 0.006312838586270617
$ws.Range("Q15").Value = 0.4567300757742221
$ws.Range("R15").Value = 4.110570681967999
$ws.Range("S15").Value = 0.005970559918285255
$ws.Range("T15").Value = 0.005970559918285254
$ws.Range("G16").Value = 1.011383666666666
$ws.Range("H16").Value = 3.034151
$ws.Range("I16").Value = 0.9457805449469654
$ws.Range("J16").Value = 0.9457805449469653
$ws.Range("M16").Value = 7.781650666666667
$ws.Range("N16").Value = 23.344952
$ws.Range("O16").Value = 0.1087809232135948
$ws.Range("P16").Value = 0.1087809232135948
$ws.Range("Q16").Value = 7.870234383972443
$ws.Range("R16").Value = 70.83210945575199
$ws.Range("S16").Value = 0.1028828808367877
$ws.Range("T16").Value = 0.1028828808367877
$ws.Range("G17").Value = 1.011383666666666
$ws.Range("H17").Value = 3.034151
$ws.Range("I17").Value = 0.9457805449469654
$ws.Range("J17").Value = 0.9457805449469653
$ws.Range("M17").Value = 3.892567333333333
$ws.Range("N17").Value = 11.677702
$ws.Range("O17").Value = 0.05441481329981927
$ws.Range("P17").Value = 0.05441481329981927
$ws.Range("Q17").Value = 3.936879022333555
$ws.Range("R17").Value = 35.43191120100199
$ws.Range("S17").Value = 0.05146447177589045
$ws.Range("T17").Value = 0.05146447177589044
$ws.Range("G18").Value = 1.011383666666666
$ws.Range("H18").Value = 3.034151
$ws.Range("I18").Value = 0.9457805449469654
$ws.Range("J18").Value = 0.9457805449469653
$ws.Range("M18").Value = 43.49559133333333
$ws.Range("N18").Value = 130.486774
$ws.Range("O18").Value = 0.6080317382054886
$ws.Range("P18").Value = 0.6080317382054886
$ws.Range("Q18").Value = 43.99073064654154
$ws.Range("R18").Value = 395.9165758188739
$ws.Range("S18").Value = 0.5750645887050376
$ws.Range("T18").Value = 0.5750645887050375
$ws.Range("G19").Value = 1.011383666666666
$ws.Range("H19").Value = 3.034151
$ws.Range("I19").Value = 0.9457805449469654
$ws.Range("J19").Value = 0.9457805449469653
$ws.Range("M19").Value = 7.007517
$ws.Range("N19").Value = 21.022551
$ws.Range("O19").Value = 0.09795918646930096
$ws.Range("P19").Value = 0.09795918646930096
$ws.Range("Q19").Value = 7.087288237688998
$ws.Range("R19").Value = 63.78559413920099
$ws.Range("S19").Value = 0.09264789276149686
$ws.Range("T19").Value = 0.09264789276149685
